# Summer 2024 Working Hours - add two new log rows (68 and 69)
# describing the PyTorch CNN work, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 68: Aug 7, 2024 (Wed), 4 hrs ---------------------------------
$ws.Cells.Item(68, 1).Value = 45511          # Date serial -> 2024-08-07
$ws.Cells.Item(68, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(68, 2).Value = "W"            # Day of week
$ws.Cells.Item(68, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(68, 3).Value = 4              # Hours

# Set the note text first so it becomes shared-string index 101
$ws.Cells.Item(68, 5).Value = "Debugging and trying to format data into model. Need to factor everything and push as ints."
$ws.Cells.Item(68, 5).WrapText = $true

# --- Row 69: Aug 8, 2024 (Thu), 8 hrs ---------------------------------
$ws.Cells.Item(69, 1).Value = 45512          # Date serial -> 2024-08-08
$ws.Cells.Item(69, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(69, 2).Value = "T"            # Day of week
$ws.Cells.Item(69, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(69, 3).Value = 8              # Hours

# The link was entered before the note text (matches original shared
# string ordering: link -> index 102, note -> index 103).
$ws.Cells.Item(69, 6).Value = "https://www.researchgate.net/figure/Model-description-of-the-hybrid-2D-CNN-Model_fig1_359461812"

$ws.Cells.Item(69, 5).Value = "Saved a factored model to use in pytorch. Reading research papers on CNNs use in grade prediction. Set up and run custom CNNs on grades. Haven't set up certain courses yet."
$ws.Cells.Item(69, 5).WrapText = $true

# --- Row heights (author resized rows while typing wrapped notes) ----
$ws.Rows.Item(68).RowHeight = 45
$ws.Rows.Item(69).RowHeight = 60

# --- Selection mirrors where the author left the cursor ---------------
$ws.Range("F68").Select() | Out-Null
